$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the two existing DV rows to "Relative performance on ..." ---
$ws.Range("B2").Value = "Relative performance on NEEAP"
$ws.Range("B3").Value = "Relative performance on NREAP "

# --- Insert two new rows for "Absolute performance" DVs right after row 3 ---
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Row 4: Absolute performance on NEEAP (EE)
$ws.Range("A4").Value = 17
$ws.Range("B4").Value = "Absolute performance on NEEAP"
$ws.Range("D4").Value = "EE_perf_abs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "EE.perfor.abs"
$ws.Range("G4").Value = "DV"
$ws.Range("J4").Value = "absolute_progress_EE"

# Row 5: Absolute performance on NREAP (REN)
$ws.Range("A5").Value = 18
$ws.Range("B5").Value = "Absolute performance on NREAP"
$ws.Range("D5").Value = "RE_perf_abs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "REN.perfor.abs"
$ws.Range("G5").Value = "DV"
$ws.Range("I5").Value = "absolute_progress_REN"

# --- Update the "Policy density" row (now row 7) : density -> number of policies ---
$ws.Range("C7").Value = "Number of policies"
$ws.Range("E7").Value = 1
$ws.Range("I7").Value = "REN_policies"
$ws.Range("J7").Value = "EE_policies"
$ws.Range("C7").Interior.Color = 65535

# --- Workbook view / window bookkeeping to match the edited file ---
$ws.Range("C11").Select() | Out-Null
